$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptos list: prices (col D) and volume-1h deltas (col E),
# plus two rows (39/40) whose coin/link also changed (Celestia <-> LidoDAOToken swap).
# Purely-numeric-looking price strings are written with a leading apostrophe so
# Excel keeps them as text (matching the original inlineStr cell type) instead of
# coercing to a number and dropping formatting like trailing zeros.

$ws.Range("D2").Value = "39.246.24"
$ws.Range("E2").Value = "  -3.28%  "
$ws.Range("D3").Value = "2.203.54"
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'295.45"
$ws.Range("E5").Value = "  -5.01%  "
$ws.Range("D6").Value = "'82.71"
$ws.Range("E6").Value = "  -4.16%  "
$ws.Range("D7").Value = "'0.509"
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  -5.62%  "
$ws.Range("D10").Value = "'0.0772"
$ws.Range("E10").Value = "  -7.86%  "
$ws.Range("D11").Value = "'28.91"
$ws.Range("E11").Value = "  -4.41%  "
$ws.Range("D12").Value = "'47.35"
$ws.Range("E12").Value = "  -10.27%  "
$ws.Range("E13").Value = "  -2.17%  "
$ws.Range("D14").Value = "2.542.49"
$ws.Range("E14").Value = "  -6.99%  "
$ws.Range("D15").Value = "'6.24"
$ws.Range("E15").Value = "  -4.47%  "
$ws.Range("D16").Value = "'14.08"
$ws.Range("E16").Value = "  -6.71%  "
$ws.Range("D17").Value = "2.195.91"
$ws.Range("E17").Value = "  -7.54%  "
$ws.Range("D18").Value = "'0.712"
$ws.Range("E18").Value = "  -6.01%  "
$ws.Range("D19").Value = "39.126.16"
$ws.Range("E19").Value = "  -3.43%  "
$ws.Range("D20").Value = "0.0₃0870"
$ws.Range("E20").Value = "  -4.43%  "
$ws.Range("D21").Value = "'5.68"
$ws.Range("E21").Value = "  -7.40%  "
$ws.Range("D22").Value = "'64.54"
$ws.Range("E22").Value = "  -5.92%  "
$ws.Range("D23").Value = "'10.27"
$ws.Range("E23").Value = "  -4.55%  "
$ws.Range("D24").Value = "'227.98"
$ws.Range("E25").Value = "  -0.10%  "
$ws.Range("E26").Value = "  -7.22%  "
$ws.Range("E27").Value = "  -1.42%  "
$ws.Range("D28").Value = "'22.49"
$ws.Range("E28").Value = "  -5.89%  "
$ws.Range("E29").Value = "  -1.68%  "
$ws.Range("D30").Value = "'9.09"
$ws.Range("E30").Value = "  -1.95%  "
$ws.Range("D31").Value = "'149.14"
$ws.Range("E31").Value = "  -3.40%  "
$ws.Range("D32").Value = "'31.63"
$ws.Range("E32").Value = "  -7.45%  "
$ws.Range("D33").Value = "'1.00"
$ws.Range("E33").Value = "  -0.06%  "
$ws.Range("E34").Value = "  -7.30%  "
$ws.Range("E35").Value = "  -5.17%  "
$ws.Range("E36").Value = "  -3.76%  "
$ws.Range("E37").Value = "  -3.93%  "
$ws.Range("D38").Value = "'0.0957"
$ws.Range("E38").Value = "  -4.61%  "
$ws.Range("B39").Value = "Celestia"
$ws.Range("C39").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D39").Value = "'15.09"
$ws.Range("E39").Value = "  -6.91%  "
$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D40").Value = "'2.62"
$ws.Range("E40").Value = "  -4.70%  "
$ws.Range("E41").Value = "  -5.34%  "
$ws.Range("D42").Value = "'3.65"
$ws.Range("E42").Value = "  -5.00%  "
$ws.Range("D43").Value = "1.901.88"
$ws.Range("E43").Value = "  -2.91%  "
$ws.Range("E44").Value = "  -4.05%  "
$ws.Range("D45").Value = "'2.02"
$ws.Range("E45").Value = "  -15.64%  "
$ws.Range("D46").Value = "'8.88"
$ws.Range("E46").Value = "  -5.13%  "
$ws.Range("D48").Value = "'15.83"
$ws.Range("E48").Value = "  -11.29%  "
$ws.Range("D49").Value = "2.408.66"
$ws.Range("E49").Value = "  -7.22%  "
$ws.Range("D50").Value = "'70.23"
$ws.Range("E50").Value = "  -2.38%  "
$ws.Range("D51").Value = "'86.34"
$ws.Range("E51").Value = "  -7.16%  "
